# Wavefolder SMD BOM - rev 2 changes
#
# Insert a new BOM line for the 1N4148 diode (D1) above the existing
# BAV99 line, and record that D13 now also uses the BAV99 part.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 4 (pushes the BAV99 row, and everything below
# it, down by one; Excel carries each row's formatting down with it).
$ws.Rows("4").Insert()

# New BOM entry: 1N4148 diode, designator D1, SOD-123 footprint, LCSC C81598.
$ws.Range("A4").Value = "1N4148"
$ws.Range("B4").Value = "D1"
$ws.Range("C4").Value = "SOD-123"
$ws.Range("D4").Value = "C81598"

# Match the red/flagged "new part" formatting already used elsewhere in
# the sheet (e.g. B2) by copying its format onto the new row.
$ws.Range("B2").Copy()
$ws.Range("A4:D4").PasteSpecial(-4122)

# BAV99 (now row 5) gains D13 as an additional designator, and picks up
# the same highlighted formatting on the Designator cell.
$ws.Range("B5").Value = "D3,D12,D11,D10,D9,D6,D5,D4,D13"
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)

$excel.CutCopyMode = $false
